# Branch 2 initial commit
# - Fill in missing "Cost" (column B) values on the "backend" sheet for rows 32-101
# - Insert a line chart on the "Revenue" sheet plotting Revenue!$E$1 (title) / $E$2:$E$10 (values)
# - Update sheet view selections/scroll positions to match the saved state

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. backend sheet: fill column B (rows 32-101) with the missing values
# ---------------------------------------------------------------------------
$backend = $wb.Worksheets.Item("backend")

$bValues = @(1000,2000,3000,4000,5000,6000,7000,8000,9000,10000,11000,12000,13000,14000,15000,16000,17000,18000,19000,20000,21000,22000,23000,24000,25000,26000,27000,28000,29000,30000,1000,2000,3000,4000,5000,6000,7000,8000,9000,10000,11000,12000,13000,14000,15000,16000,17000,18000,19000,20000,21000,22000,23000,24000,25000,26000,27000,28000,29000,30000,31000,32000,33000,34000,35000,36000,37000,38000,39000,40000)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = 32 + $i
    $backend.Cells.Item($row, 2).Value = $bValues[$i]
}

# ---------------------------------------------------------------------------
# 2. Revenue sheet: add a line chart sourced from Revenue!$E$1:$E$10
# ---------------------------------------------------------------------------
$revenue = $wb.Worksheets.Item("Revenue")
$revenue.Activate()

$chartObj = $revenue.Shapes.AddChart2(-1, 4)
$chartObj.Name = "Chart 2"

$chart = $chartObj.Chart
$chart.SetSourceData($revenue.Range("E1:E10"))
$chart.ChartType = 4

# ---------------------------------------------------------------------------
# 3. Restore the saved selections / scroll positions for each sheet
# ---------------------------------------------------------------------------
$backend.Activate()
$backend.Range("F107").Select()

$revenue.Activate()
$revenue.Range("E1:E10").Select()
